$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "| Top Secret-SCI" -> "| San Diego, Ca"
# The run holding this text sits right after a run that is a lone
# space (" "); a plain Range.Text assignment causes the engine to
# auto-merge that neighboring run into this one (since they share
# identical formatting). To keep that neighbor run intact, we briefly
# toggle Bold on it to break the formatting match, do the text swap,
# then toggle Bold back off.
# ---------------------------------------------------------------------
$marker1 = "| Top Secret-SCI"
$t = $d.Content.Text
$idx = $t.IndexOf($marker1)

$rBlock1 = $d.Range($idx - 1, $idx)
$rBlock1.Font.Bold = 1

$t = $d.Content.Text
$idx = $t.IndexOf($marker1)
$r1 = $d.Range($idx, $idx + $marker1.Length)
$r1.Text = "| San Diego, Ca"

$t = $d.Content.Text
$idx = $t.IndexOf("| San Diego, Ca")
$rBlock1b = $d.Range($idx - 1, $idx)
$rBlock1b.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2: clean up the Summary paragraph -- remove the gramStart/
# gramEnd proofErr-error markers and merge the runs they used to
# separate, while leaving the "over " run (which sits in between and
# must stay a distinct run) untouched. The engine merges a run with
# its same-formatted neighbors whenever its text is assigned via
# Range.Text, so "over " is temporarily bolded to act as a formatting
# break/blocker while the two flanking segments are rebuilt, then
# un-bolded again afterward.
#
# Because Range.Text is a no-op when the assigned string already
# equals the current text (even if it spans runs that still need to
# be merged), each segment is first swapped to a placeholder string
# and then immediately back to the real text, which forces the
# engine to actually rebuild/merge the runs.
# ---------------------------------------------------------------------
$overMarker = "over "
$t = $d.Content.Text
$overIdx = $t.IndexOf($overMarker)
$rOver = $d.Range($overIdx, $overIdx + $overMarker.Length)
$rOver.Font.Bold = 1

# Left segment: "Results-driven" + " Management professional with "
$t = $d.Content.Text
$leftStart = $t.IndexOf("Results-driven")
$overIdx = $t.IndexOf($overMarker)
$rLeft = $d.Range($leftStart, $overIdx)
$rLeft.Text = "PLACEHOLDER_LEFT_SEGMENT"

$t = $d.Content.Text
$phIdx = $t.IndexOf("PLACEHOLDER_LEFT_SEGMENT")
$rLeft2 = $d.Range($phIdx, $phIdx + "PLACEHOLDER_LEFT_SEGMENT".Length)
$rLeft2.Text = "Results-driven Management professional with "

# Right segment: "10 years ... " + "Track record" + " of improving ... objectives."
$t = $d.Content.Text
$overIdx = $t.IndexOf($overMarker)
$overEnd = $overIdx + $overMarker.Length
$endMarker = "achieving defined objectives."
$endIdx = $t.IndexOf($endMarker) + $endMarker.Length
$rRight = $d.Range($overEnd, $endIdx)
$rRight.Text = "PLACEHOLDER_RIGHT_SEGMENT"

$t = $d.Content.Text
$phIdx2 = $t.IndexOf("PLACEHOLDER_RIGHT_SEGMENT")
$rRight2 = $d.Range($phIdx2, $phIdx2 + "PLACEHOLDER_RIGHT_SEGMENT".Length)
$rightText = "10 years of experience in high-risk operations, organizational development and team building within diverse industries. Highly adept in planning, coordinating, and executing successful production strategies. Track record of improving operational stability, efficiency, and profitability. Strong collaborator with senior stakeholders, effectively prioritizing activities, translating business requirements into solutions and achieving defined objectives."
$rRight2.Text = $rightText

# Un-bold "over " again so it is left exactly as it started.
$t = $d.Content.Text
$overIdx = $t.IndexOf($overMarker)
$rOverB = $d.Range($overIdx, $overIdx + $overMarker.Length)
$rOverB.Font.Bold = 0
